# Add data for 2024-08-01
# Updates 2024 (column K) crime-count figures across the citywide, by-neighborhood,
# and per-neighborhood worksheets (plus a few 2021/column-H corrections) to reflect
# the latest August 2024 data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4629
$ws.Range("K3").Value = 4755
$ws.Range("H4").Value = 1739
$ws.Range("K4").Value = 968
$ws.Range("K6").Value = 5356
$ws.Range("H7").Value = 26052
$ws.Range("K7").Value = 16051

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 299
$ws.Range("K3").Value = 325
$ws.Range("K4").Value = 62
$ws.Range("K6").Value = 365
$ws.Range("K7").Value = 1080

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 118
$ws.Range("K7").Value = 341

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 185
$ws.Range("K3").Value = 252
$ws.Range("K7").Value = 673

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 96
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 152
$ws.Range("K3").Value = 181
$ws.Range("K7").Value = 542

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 114
$ws.Range("K3").Value = 91
$ws.Range("K7").Value = 369

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 142
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 474
$ws.Range("K8").Value = 1080
$ws.Range("K10").Value = 88
$ws.Range("K11").Value = 315
$ws.Range("K14").Value = 89
$ws.Range("K18").Value = 110
$ws.Range("K19").Value = 480
$ws.Range("K20").Value = 371
$ws.Range("K21").Value = 47
$ws.Range("K23").Value = 166
$ws.Range("K24").Value = 48
$ws.Range("K27").Value = 144
$ws.Range("K29").Value = 853
$ws.Range("K33").Value = 673
$ws.Range("K34").Value = 82
$ws.Range("K35").Value = 23
$ws.Range("K36").Value = 203
$ws.Range("K37").Value = 542
$ws.Range("K42").Value = 598
$ws.Range("K44").Value = 141
$ws.Range("K48").Value = 202
$ws.Range("K51").Value = 202
$ws.Range("K52").Value = 422
$ws.Range("K54").Value = 305
$ws.Range("H63").Value = 290
$ws.Range("K63").Value = 49
$ws.Range("K64").Value = 96
$ws.Range("K65").Value = 369
$ws.Range("K67").Value = 620
$ws.Range("K71").Value = 52
$ws.Range("K72").Value = 75
$ws.Range("K77").Value = 118
$ws.Range("K78").Value = 188
$ws.Range("K79").Value = 397
$ws.Range("K83").Value = 341
$ws.Range("K85").Value = 731
$ws.Range("K86").Value = 106
$ws.Range("K88").Value = 184
$ws.Range("K89").Value = 229
$ws.Range("K90").Value = 146
$ws.Range("K91").Value = 173
$ws.Range("K94").Value = 205
$ws.Range("K95").Value = 281
$ws.Range("K99").Value = 269
$ws.Range("K100").Value = 29
$ws.Range("H101").Value = 26052
$ws.Range("K101").Value = 16051

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 180
$ws.Range("K3").Value = 215
$ws.Range("K6").Value = 181
$ws.Range("K7").Value = 620

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 304
$ws.Range("K6").Value = 239
$ws.Range("K7").Value = 853

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 480

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 121

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 161
$ws.Range("K6").Value = 225
$ws.Range("K7").Value = 598

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 83
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 129
$ws.Range("K7").Value = 397

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 117
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 371

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 80
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 474

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 36
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 82
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 73
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 43
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 254
$ws.Range("K3").Value = 245
$ws.Range("K4").Value = 42
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 731

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 11
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 158
$ws.Range("K7").Value = 422
